# 9th Stab- Cosmetic Changes
# Insert two new "report date" columns (Jun_17, Jun_15) in front of the
# existing date columns, pushing the older Jun_13 / Jun_10 columns to the
# right. New cells default to "UN" (unchanged / no rating action yet),
# mirroring the rest of that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before the current column B (the first date
# column). This shifts the old B ("Jun_13") -> D and old C ("Jun_10") -> E.
$ws.Range("B1:C1").EntireColumn.Insert()

# New header cells for the two freshly inserted date columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the new columns' data rows with the default "UN" rating used
# throughout the rest of the sheet.
$ws.Range("B2:C27").Value = "UN"
